$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 383
$ws.Range("I2").Value = 1006
$ws.Range("J2").Value = 4188
$ws.Range("K2").Value = 33
$ws.Range("L2").Value = 1158
$ws.Range("M2").Value = 61
$ws.Range("N2").Value = 732
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 21
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 53
$ws.Range("S2").Value = 452
$ws.Range("T2").Value = 743
$ws.Range("U2").Value = 56
$ws.Range("V2").Value = 6565
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 6609
$ws.Range("Y2").Value = 14
$ws.Range("Z2").Value = 110
$ws.Range("AA2").Value = 41
